# Apply the transaction-report update:
#  - extend the "Transaction Dates" range start (07/03 -> 06/23)
#  - bump the "Report Generated on" timestamp
#  - refresh Balance / Total Expenses to reflect a new expense row
#  - fix the mislabeled Category on the Salary (Income) row
#  - append a new expense transaction row (2024-06-23, House and Garden)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary / header block (column A, rows 2-5) ---
$ws.Range("A2").Value = "Transaction Dates: " + [char]10 + "06/23/2024" + [char]0x2014 + "10/15/2024"
$ws.Range("A3").Value = "Report Generated on:" + [char]10 + "09/02/2024 13:53:40"
$ws.Range("A4").Value = "Balance:" + [char]10 + "1679.01"
$ws.Range("A5").Value = "Total Expenses:" + [char]10 + "7520.99"

# --- Fix mislabeled category on the Salary / Income row (row 14) ---
$ws.Range("C14").Value = "Salary"

# --- New expense transaction row appended as row 17 ---
# Force the date-looking / number-looking text into real text cells (like
# the rest of column A / D) by pre-formatting as Text, then re-apply the
# plain formatting from the row above so no stray number format lingers.
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "2024-06-23"
$ws.Range("B17").Value = "Expense"
$ws.Range("C17").Value = "House and Garden"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "150.99"
$ws.Range("E17").Value = "Garden hose and watering can for my garden"

$ws.Range("A16:E16").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
